$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 1, 6.201049113329182)
    3 = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 1, 6.201049113329182)
    4 = @(1.459612070389937, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 0, 3.781711156805759)
    5 = @(3.230985683306322, 1.667794583268128, 3.900430680208489, 8.660232485948974, 1, 17.45944343273191)
    6 = @(1.459612070389937, 1.667794583268128, 26.21740644021617, 8.660232485948974, 0, 38.00504557982321)
    7 = @(0.003994804209775715, 0.002777888934908601, 0.8054896365839992, 0.496779210170732, 1, 1.309041539899416)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals[0]
    $ws.Range("C$row").Value = $vals[1]
    $ws.Range("D$row").Value = $vals[2]
    $ws.Range("E$row").Value = $vals[3]
    $ws.Range("F$row").Value = $vals[4]
    $ws.Range("G$row").Value = $vals[5]
}
